$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 10 (current row 10 "Especial" / 44273 and
# row 11 "Primera" / 44273 shift down to rows 12/13, and so on through the
# old row 27 which becomes row 29).
$ws.Range("A10:T11").EntireRow.Insert()

# New row 10: Especial quality, Región de O'Higgins, week of 2022-04-11
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44662
$ws.Cells.Item(10, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = "Fruta"
$ws.Cells.Item(10, 7).Value = 100104
$ws.Cells.Item(10, 8).Value = "Frutos de pepita"
$ws.Cells.Item(10, 9).Value = 100104003
$ws.Cells.Item(10, 10).Value = "Membrillo"
$ws.Cells.Item(10, 11).Value = "Champion"
$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 20
$ws.Cells.Item(10, 14).Value = 330000
$ws.Cells.Item(10, 15).Value = 340000
$ws.Cells.Item(10, 16).Value = 335000
$ws.Cells.Item(10, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(10, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(10, 19).Value = 744
$ws.Cells.Item(10, 20).Value = 450

# New row 11: Primera quality, Región de O'Higgins, week of 2022-04-11
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44662
$ws.Cells.Item(11, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100104
$ws.Cells.Item(11, 8).Value = "Frutos de pepita"
$ws.Cells.Item(11, 9).Value = 100104003
$ws.Cells.Item(11, 10).Value = "Membrillo"
$ws.Cells.Item(11, 11).Value = "Champion"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 20
$ws.Cells.Item(11, 14).Value = 300000
$ws.Cells.Item(11, 15).Value = 310000
$ws.Cells.Item(11, 16).Value = 305000
$ws.Cells.Item(11, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(11, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(11, 19).Value = 678
$ws.Cells.Item(11, 20).Value = 450
